$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to hold the given string as plain text, preserving
# the original (default) cell style, so values like "1.00" or "0.401"
# are not silently coerced into numeric doubles by Excel.
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$sub3 = [string][char]0x2083

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "68.393.48"
$ws.Range("E2").Value = "  -3.05%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "3.708.07"
$ws.Range("E3").Value = "  -3.72%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "602.67"
$ws.Range("E5").Value = "  +1.89%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "181.01"
$ws.Range("E6").Value = "  +8.76%  "

# Row 7 - LidoStakedEther
Set-TextValue $ws.Range("D7") "3.696.15"
$ws.Range("E7").Value = "  -3.79%  "

# Row 8 - XRP
Set-TextValue $ws.Range("D8") "0.632"
$ws.Range("E8").Value = "  -5.83%  "

# Row 9 - USDC
Set-TextValue $ws.Range("D9") "1.00"
$ws.Range("E9").Value = "  -0.03%  "

# Row 10 - Cardano
Set-TextValue $ws.Range("D10") "0.718"
$ws.Range("E10").Value = "  -3.98%  "

# Row 11 - Dogecoin
Set-TextValue $ws.Range("D11") "0.163"
$ws.Range("E11").Value = "  -6.71%  "

# Row 12 - Avalanche
Set-TextValue $ws.Range("D12") "56.41"
$ws.Range("E12").Value = "  +6.30%  "

# Row 13 - ShibaInu
Set-TextValue $ws.Range("D13") "0.0000293"
$ws.Range("E13").Value = "  -8.21%  "

# Row 14 - Polkadot
Set-TextValue $ws.Range("D14") "10.44"
$ws.Range("E14").Value = "  -7.98%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D15") "4.290.14"
$ws.Range("E15").Value = "  -4.06%  "

# Row 16 - WrappedEther
Set-TextValue $ws.Range("D16") "3.705.14"
$ws.Range("E16").Value = "  -4.05%  "

# Row 17 - Chainlink
Set-TextValue $ws.Range("D17") "19.37"
$ws.Range("E17").Value = "  -8.20%  "

# Row 18 - TRON
$ws.Range("E18").Value = "  -2.16%  "

# Row 19 - Uniswap
Set-TextValue $ws.Range("D19") "12.90"
$ws.Range("E19").Value = "  -6.28%  "

# Row 20 - Polygon
Set-TextValue $ws.Range("D20") "1.13"
$ws.Range("E20").Value = "  -6.62%  "

# Row 21 - WrappedBTC
Set-TextValue $ws.Range("D21") "68.232.34"
$ws.Range("E21").Value = "  -3.27%  "

# Row 22 - BitcoinCash
Set-TextValue $ws.Range("D22") "409.78"
$ws.Range("E22").Value = "  -6.16%  "

# Row 23 - PancakeSwap
Set-TextValue $ws.Range("D23") "4.64"
$ws.Range("E23").Value = "  -1.35%  "

# Row 24 - Litecoin
Set-TextValue $ws.Range("D24") "89.35"
$ws.Range("E24").Value = "  -4.75%  "

# Row 25 - ImmutableX
Set-TextValue $ws.Range("D25") "3.03"
$ws.Range("E25").Value = "  -6.79%  "

# Row 26 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D26") "12.88"
$ws.Range("E26").Value = "  -6.86%  "

# Row 27 - RenderToken
Set-TextValue $ws.Range("D27") "10.86"
$ws.Range("E27").Value = "  -2.72%  "

# Row 28 - Toncoin
Set-TextValue $ws.Range("D28") "3.90"
$ws.Range("E28").Value = "  -3.63%  "

# Row 29 - LEO
Set-TextValue $ws.Range("D29") "6.06"
$ws.Range("E29").Value = "  +1.97%  "

# Row 30 - Filecoin
Set-TextValue $ws.Range("D30") "9.49"
$ws.Range("E30").Value = "  -8.62%  "

# Row 31 - EthereumClassic
Set-TextValue $ws.Range("D31") "32.95"
$ws.Range("E31").Value = "  -5.97%  "

# Row 32 - NEARProtocol
Set-TextValue $ws.Range("D32") "7.31"
$ws.Range("E32").Value = "  -10.25%  "

# Row 33 - Cosmos
Set-TextValue $ws.Range("D33") "12.52"
$ws.Range("E33").Value = "  -7.06%  "

# Row 34 - Hedera
Set-TextValue $ws.Range("D34") "0.118"
$ws.Range("E34").Value = "  -5.70%  "

# Row 35 - InjectiveProtocol
Set-TextValue $ws.Range("D35") "43.91"
$ws.Range("E35").Value = "  -8.86%  "

# Row 36 - OKB
Set-TextValue $ws.Range("D36") "64.42"
$ws.Range("E36").Value = "  -7.29%  "

# Row 37 - Bittensor
Set-TextValue $ws.Range("D37") "601.77"
$ws.Range("E37").Value = "  -4.35%  "

# Row 38 - PEPE
$pepePrice = "0.0" + $sub3 + "0892"
Set-TextValue $ws.Range("D38") $pepePrice
$ws.Range("E38").Value = "  -8.37%  "

# Row 39 - now TheGraph (was Dai)
$ws.Range("B39").Value = "TheGraph"
$ws.Range("C39").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue $ws.Range("D39") "0.401"
$ws.Range("E39").Value = "  -4.87%  "

# Row 40 - now Dai (was TheGraph)
$ws.Range("B40").Value = "Dai"
$ws.Range("C40").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D40") "1.00"
$ws.Range("E40").Value = "  +0.25%  "

# Row 41 - FirstDigitalUSD
$ws.Range("E41").Value = "  -0.05%  "

# Row 42 - Kaspa
Set-TextValue $ws.Range("D42") "0.137"
$ws.Range("E42").Value = "  -5.76%  "

# Row 43 - Fetch.AI
Set-TextValue $ws.Range("D43") "2.76"
$ws.Range("E43").Value = "  +2.81%  "

# Row 44 - ThetaToken
Set-TextValue $ws.Range("D44") "3.04"
$ws.Range("E44").Value = "  -7.74%  "

# Row 45 - VeChain
Set-TextValue $ws.Range("D45") "0.0439"
$ws.Range("E45").Value = "  -6.34%  "

# Row 46 - dogwifhat
Set-TextValue $ws.Range("D46") "2.89"
$ws.Range("E46").Value = "  -12.00%  "

# Row 47 - THORChain
Set-TextValue $ws.Range("D47") "9.24"
$ws.Range("E47").Value = "  -7.72%  "

# Row 48 - WEMIXToken
Set-TextValue $ws.Range("D48") "2.73"
$ws.Range("E48").Value = "  -3.68%  "

# Row 49 - Stellar
Set-TextValue $ws.Range("D49") "0.135"
$ws.Range("E49").Value = "  -5.96%  "

# Row 50 - Maker
Set-TextValue $ws.Range("D50") "2.775.58"
$ws.Range("E50").Value = "  -2.01%  "

# Row 51 - ApeXProtocol
$ws.Range("E51").Value = "  -2.24%  "
